# The commit swaps the contents of ppt/theme/theme1.xml (the slide
# master's theme, originally the "Integral" theme) and
# ppt/theme/theme2.xml (the notes master's theme, originally the
# default "Office Theme"): after the edit, theme1.xml carries the
# "Office Theme" palette and theme2.xml carries the "Integral" palette.
#
# The font scheme (<a:fontScheme>) and format scheme (<a:fmtScheme>)
# are byte-for-byte identical between the two themes, so the only real
# difference is the color scheme (<a:clrScheme>) - the 12 theme
# colors. We reproduce that by rewriting the slide master's theme
# colors to the "Office Theme" palette via the PowerPoint object
# model.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Index order (verified): 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2
# 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
# RGB() packs as 0x00BBGGRR to match the COM OLE_COLOR layout used by
# RGBColor.RGB.
function Pack-RGB($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$officeTheme = @{
    1  = Pack-RGB 0x00 0x00 0x00   # dk1
    2  = Pack-RGB 0xFF 0xFF 0xFF   # lt1
    3  = Pack-RGB 0x44 0x54 0x6A   # dk2
    4  = Pack-RGB 0xE7 0xE6 0xE6   # lt2
    5  = Pack-RGB 0x5B 0x9B 0xD5   # accent1
    6  = Pack-RGB 0xED 0x7D 0x31   # accent2
    7  = Pack-RGB 0xA5 0xA5 0xA5   # accent3
    8  = Pack-RGB 0xFF 0xC0 0x00   # accent4
    9  = Pack-RGB 0x44 0x72 0xC4   # accent5
    10 = Pack-RGB 0x70 0xAD 0x47   # accent6
    11 = Pack-RGB 0x05 0x63 0xC1   # hlink
    12 = Pack-RGB 0x95 0x4F 0x72   # folHlink
}

for ($i = 1; $i -le 12; $i++) {
    $colors.Item($i).RGB = $officeTheme[$i]
}
